$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 66
$ws.Range("I2").Value = 212
$ws.Range("J2").Value = 809
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 235
$ws.Range("M2").Value = 11
$ws.Range("N2").Value = 135
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 5
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 9
$ws.Range("S2").Value = 86
$ws.Range("T2").Value = 126
$ws.Range("U2").Value = 9
$ws.Range("V2").Value = 1190
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 1239
$ws.Range("Y2").Value = 0
$ws.Range("Z2").Value = 21
$ws.Range("AA2").Value = 7
